# Insert a new weekly record at row 544 of "Fruta, Feria Lagunitas de Puerto
# Montt - Pomelo", pushing the existing rows 544-671 down to 545-672 (new
# dimension A1:T672). The new row carries the same market/product metadata
# as the rest of the block (columns A,B,C,E,F,G,H,I,J,K,Q,R,T) together with
# its own date/quality/volume/price/weighted-price data (columns D,L,M,N,O,P,S).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Shift rows 544:671 down by one, creating a blank row 544.
$ws.Rows.Item(544).Insert()

# Fill in the new row 544 with the new observation.
$ws.Cells.Item(544, 1).Value  = 4
$ws.Cells.Item(544, 2).Value  = "Feria Lagunitas de Puerto Montt"
$ws.Cells.Item(544, 3).Value  = "Los Lagos"
$ws.Cells.Item(544, 4).Value  = 45173
$ws.Cells.Item(544, 5).Value  = 10
$ws.Cells.Item(544, 6).Value  = "Fruta"
$ws.Cells.Item(544, 7).Value  = 100102
$ws.Cells.Item(544, 8).Value  = "Cítricos"
$ws.Cells.Item(544, 9).Value  = 100102006
$ws.Cells.Item(544, 10).Value = "Pomelo"
$ws.Cells.Item(544, 11).Value = "Start Ruby"
$ws.Cells.Item(544, 12).Value = "Primera"
$ws.Cells.Item(544, 13).Value = 80
$ws.Cells.Item(544, 14).Value = 14000
$ws.Cells.Item(544, 15).Value = 15000
$ws.Cells.Item(544, 16).Value = 14500
$ws.Cells.Item(544, 17).Value = "$/caja 14 kilos empedrada"
$ws.Cells.Item(544, 18).Value = "Región de O'Higgins"
$ws.Cells.Item(544, 19).Value = 1036
$ws.Cells.Item(544, 20).Value = 14
